# Implements: "Implemented getting number of lines for methods and classes."
#
# 1) Adds two new sheets: classNumberOfLines, methodNumberOfLines
# 2) Reorders some field rows within the existing classFields sheet
#    (reflects the underlying Java source field order changing when the
#    code was rescanned to also compute line counts).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1a. classNumberOfLines sheet
# ---------------------------------------------------------------------
$classSheet = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$classSheet.Name = "classNumberOfLines"

$classSheet.Range("A1").Value = "Class Name"
$classSheet.Range("B1").Value = "Number of Lines"

$classRows = @(
    @("com.macro.mall.search.domain.EsProductAttributeValue", "39"),
    @("com.macro.mall.search.controller.EsProductController", "39"),
    @("com.macro.mall.search.service.EsProductService", "11"),
    @("com.macro.mall.search.config.SwaggerConfig", "9"),
    @("com.macro.mall.search.domain.EsProduct", "130"),
    @("com.macro.mall.search.service.impl.EsProductServiceImpl", "187"),
    @("com.macro.mall.search.domain.EsProductRelatedInfo`$ProductAttr", "24"),
    @("com.macro.mall.search.domain.EsProductRelatedInfo", "48"),
    @("com.macro.mall.search.dao.EsProductDao", "4"),
    @("com.macro.mall.search.repository.EsProductRepository", "4"),
    @("com.macro.mall.search.MallSearchApplication", "6"),
    @("com.macro.mall.search.config.MyBatisConfig", "3")
)

$r = 2
foreach ($row in $classRows) {
    $classSheet.Cells.Item($r, 1).Value = $row[0]
    $classSheet.Cells.Item($r, 2).Formula = "'" + $row[1]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 1b. methodNumberOfLines sheet
# ---------------------------------------------------------------------
$methodSheet = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$methodSheet.Name = "methodNumberOfLines"

$methodSheet.Range("A1").Value = "Class Name"
$methodSheet.Range("B1").Value = "Method Signature"
$methodSheet.Range("C1").Value = "Number of Lines"

$methodRows = @(
    @("com.macro.mall.search.domain.EsProductAttributeValue", "getId()", "3"),
    @("com.macro.mall.search.domain.EsProductAttributeValue", "setId(java.lang.Long)", "3"),
    @("com.macro.mall.search.domain.EsProductAttributeValue", "getProductAttributeId()", "3"),
    @("com.macro.mall.search.domain.EsProductAttributeValue", "setProductAttributeId(java.lang.Long)", "3"),
    @("com.macro.mall.search.domain.EsProductAttributeValue", "getValue()", "3"),
    @("com.macro.mall.search.domain.EsProductAttributeValue", "setValue(java.lang.String)", "3"),
    @("com.macro.mall.search.domain.EsProductAttributeValue", "getType()", "3"),
    @("com.macro.mall.search.domain.EsProductAttributeValue", "setType(java.lang.Integer)", "3"),
    @("com.macro.mall.search.domain.EsProductAttributeValue", "getName()", "3"),
    @("com.macro.mall.search.domain.EsProductAttributeValue", "setName(java.lang.String)", "3"),
    @("com.macro.mall.search.controller.EsProductController", "importAllList()", "4"),
    @("com.macro.mall.search.controller.EsProductController", "delete(java.lang.Long)", "4"),
    @("com.macro.mall.search.controller.EsProductController", "delete(java.util.List)", "4"),
    @("com.macro.mall.search.controller.EsProductController", "create(java.lang.Long)", "5"),
    @("com.macro.mall.search.controller.EsProductController", "search(java.lang.String, java.lang.Integer, java.lang.Integer)", "4"),
    @("com.macro.mall.search.controller.EsProductController", "search(java.lang.String, java.lang.Long, java.lang.Long, java.lang.Integer, java.lang.Integer, java.lang.Integer)", "4"),
    @("com.macro.mall.search.controller.EsProductController", "recommend(java.lang.Long, java.lang.Integer, java.lang.Integer)", "4"),
    @("com.macro.mall.search.controller.EsProductController", "searchRelatedInfo(java.lang.String)", "4"),
    @("com.macro.mall.search.service.EsProductService", "importAll()", "1"),
    @("com.macro.mall.search.service.EsProductService", "delete(java.lang.Long)", "1"),
    @("com.macro.mall.search.service.EsProductService", "create(java.lang.Long)", "1"),
    @("com.macro.mall.search.service.EsProductService", "delete(java.util.List)", "1"),
    @("com.macro.mall.search.service.EsProductService", "search(java.lang.String, java.lang.Integer, java.lang.Integer)", "1"),
    @("com.macro.mall.search.service.EsProductService", "search(java.lang.String, java.lang.Long, java.lang.Long, java.lang.Integer, java.lang.Integer, java.lang.Integer)", "1"),
    @("com.macro.mall.search.service.EsProductService", "recommend(java.lang.Long, java.lang.Integer, java.lang.Integer)", "1"),
    @("com.macro.mall.search.service.EsProductService", "searchRelatedInfo(java.lang.String)", "1"),
    @("com.macro.mall.search.config.SwaggerConfig", "swaggerProperties()", "3"),
    @("com.macro.mall.search.config.SwaggerConfig", "springfoxHandlerProviderBeanPostProcessor()", "3"),
    @("com.macro.mall.search.domain.EsProduct", "getId()", "3"),
    @("com.macro.mall.search.domain.EsProduct", "setId(java.lang.Long)", "3"),
    @("com.macro.mall.search.domain.EsProduct", "getProductSn()", "3"),
    @("com.macro.mall.search.domain.EsProduct", "setProductSn(java.lang.String)", "3"),
    @("com.macro.mall.search.domain.EsProduct", "getBrandId()", "3"),
    @("com.macro.mall.search.domain.EsProduct", "setBrandId(java.lang.Long)", "3"),
    @("com.macro.mall.search.domain.EsProduct", "getBrandName()", "3"),
    @("com.macro.mall.search.domain.EsProduct", "setBrandName(java.lang.String)", "3"),
    @("com.macro.mall.search.domain.EsProduct", "getProductCategoryId()", "3"),
    @("com.macro.mall.search.domain.EsProduct", "setProductCategoryId(java.lang.Long)", "3"),
    @("com.macro.mall.search.domain.EsProduct", "getProductCategoryName()", "3"),
    @("com.macro.mall.search.domain.EsProduct", "setProductCategoryName(java.lang.String)", "3"),
    @("com.macro.mall.search.domain.EsProduct", "getPic()", "3"),
    @("com.macro.mall.search.domain.EsProduct", "setPic(java.lang.String)", "3"),
    @("com.macro.mall.search.domain.EsProduct", "getName()", "3"),
    @("com.macro.mall.search.domain.EsProduct", "setName(java.lang.String)", "3"),
    @("com.macro.mall.search.domain.EsProduct", "getSubTitle()", "3"),
    @("com.macro.mall.search.domain.EsProduct", "setSubTitle(java.lang.String)", "3"),
    @("com.macro.mall.search.domain.EsProduct", "getPrice()", "3"),
    @("com.macro.mall.search.domain.EsProduct", "setPrice(java.math.BigDecimal)", "3"),
    @("com.macro.mall.search.domain.EsProduct", "getSale()", "3"),
    @("com.macro.mall.search.domain.EsProduct", "setSale(java.lang.Integer)", "3"),
    @("com.macro.mall.search.domain.EsProduct", "getNewStatus()", "3"),
    @("com.macro.mall.search.domain.EsProduct", "setNewStatus(java.lang.Integer)", "3"),
    @("com.macro.mall.search.domain.EsProduct", "getRecommandStatus()", "3"),
    @("com.macro.mall.search.domain.EsProduct", "setRecommandStatus(java.lang.Integer)", "3"),
    @("com.macro.mall.search.domain.EsProduct", "getStock()", "3"),
    @("com.macro.mall.search.domain.EsProduct", "setStock(java.lang.Integer)", "3"),
    @("com.macro.mall.search.domain.EsProduct", "getPromotionType()", "3"),
    @("com.macro.mall.search.domain.EsProduct", "setPromotionType(java.lang.Integer)", "3"),
    @("com.macro.mall.search.domain.EsProduct", "getSort()", "3"),
    @("com.macro.mall.search.domain.EsProduct", "setSort(java.lang.Integer)", "3"),
    @("com.macro.mall.search.domain.EsProduct", "getAttrValueList()", "3"),
    @("com.macro.mall.search.domain.EsProduct", "setAttrValueList(java.util.List)", "3"),
    @("com.macro.mall.search.domain.EsProduct", "getKeywords()", "3"),
    @("com.macro.mall.search.domain.EsProduct", "setKeywords(java.lang.String)", "3"),
    @("com.macro.mall.search.service.impl.EsProductServiceImpl", "importAll()", "9"),
    @("com.macro.mall.search.service.impl.EsProductServiceImpl", "delete(java.lang.Long)", "3"),
    @("com.macro.mall.search.service.impl.EsProductServiceImpl", "create(java.lang.Long)", "7"),
    @("com.macro.mall.search.service.impl.EsProductServiceImpl", "delete(java.util.List)", "8"),
    @("com.macro.mall.search.service.impl.EsProductServiceImpl", "search(java.lang.String, java.lang.Integer, java.lang.Integer)", "4"),
    @("com.macro.mall.search.service.impl.EsProductServiceImpl", "search(java.lang.String, java.lang.Long, java.lang.Long, java.lang.Integer, java.lang.Integer, java.lang.Integer)", "9"),
    @("com.macro.mall.search.service.impl.EsProductServiceImpl", "recommend(java.lang.Long, java.lang.Integer, java.lang.Integer)", "6"),
    @("com.macro.mall.search.service.impl.EsProductServiceImpl", "searchRelatedInfo(java.lang.String)", "5"),
    @("com.macro.mall.search.service.impl.EsProductServiceImpl", "convertProductRelatedInfo(org.springframework.data.elasticsearch.core.SearchHits)", "8"),
    @("com.macro.mall.search.domain.EsProductRelatedInfo`$ProductAttr", "getAttrId()", "3"),
    @("com.macro.mall.search.domain.EsProductRelatedInfo`$ProductAttr", "setAttrId(java.lang.Long)", "3"),
    @("com.macro.mall.search.domain.EsProductRelatedInfo`$ProductAttr", "getAttrValues()", "3"),
    @("com.macro.mall.search.domain.EsProductRelatedInfo`$ProductAttr", "setAttrValues(java.util.List)", "3"),
    @("com.macro.mall.search.domain.EsProductRelatedInfo`$ProductAttr", "getAttrName()", "3"),
    @("com.macro.mall.search.domain.EsProductRelatedInfo`$ProductAttr", "setAttrName(java.lang.String)", "3"),
    @("com.macro.mall.search.domain.EsProductRelatedInfo", "getBrandNames()", "3"),
    @("com.macro.mall.search.domain.EsProductRelatedInfo", "setBrandNames(java.util.List)", "3"),
    @("com.macro.mall.search.domain.EsProductRelatedInfo", "getProductCategoryNames()", "3"),
    @("com.macro.mall.search.domain.EsProductRelatedInfo", "setProductCategoryNames(java.util.List)", "3"),
    @("com.macro.mall.search.domain.EsProductRelatedInfo", "getProductAttrs()", "3"),
    @("com.macro.mall.search.domain.EsProductRelatedInfo", "setProductAttrs(java.util.List)", "3"),
    @("com.macro.mall.search.dao.EsProductDao", "getAllEsProductList(java.lang.Long)", "1"),
    @("com.macro.mall.search.repository.EsProductRepository", "findByNameOrSubTitleOrKeywords(java.lang.String, java.lang.String, java.lang.String, org.springframework.data.domain.Pageable)", "1"),
    @("com.macro.mall.search.MallSearchApplication", "main(java.lang.String[])", "3")
)

$r = 2
foreach ($row in $methodRows) {
    $methodSheet.Cells.Item($r, 1).Value = $row[0]
    $methodSheet.Cells.Item($r, 2).Value = $row[1]
    $methodSheet.Cells.Item($r, 3).Formula = "'" + $row[2]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2. classFields sheet: field order within classes changed
#    (Class Name / Field Modifier columns are unaffected; only a
#    subset of rows need their Field Name + Field Type updated.)
# ---------------------------------------------------------------------
$fieldsSheet = $wb.Worksheets.Item("classFields")

$fieldChanges = @(
    @(2,  "productAttributeId", "java.lang.Long"),
    @(4,  "name", "java.lang.String"),
    @(5,  "id", "java.lang.Long"),
    @(6,  "value", "java.lang.String"),
    @(7,  "type", "java.lang.Integer"),
    @(9,  "keywords", "java.lang.String"),
    @(10, "newStatus", "java.lang.Integer"),
    @(11, "brandName", "java.lang.String"),
    @(12, "attrValueList", "java.util.List"),
    @(13, "name", "java.lang.String"),
    @(14, "serialVersionUID", "long"),
    @(15, "brandId", "java.lang.Long"),
    @(16, "sort", "java.lang.Integer"),
    @(17, "subTitle", "java.lang.String"),
    @(18, "price", "java.math.BigDecimal"),
    @(19, "sale", "java.lang.Integer"),
    @(20, "promotionType", "java.lang.Integer"),
    @(21, "recommandStatus", "java.lang.Integer"),
    @(22, "id", "java.lang.Long"),
    @(23, "productSn", "java.lang.String"),
    @(24, "productCategoryId", "java.lang.Long"),
    @(25, "stock", "java.lang.Integer"),
    @(26, "pic", "java.lang.String"),
    @(27, "productCategoryName", "java.lang.String"),
    @(28, "productDao", "com.macro.mall.search.dao.EsProductDao"),
    @(31, "elasticsearchRestTemplate", "org.springframework.data.elasticsearch.core.ElasticsearchRestTemplate")
)

foreach ($chg in $fieldChanges) {
    $rowIdx = $chg[0]
    $fieldsSheet.Cells.Item($rowIdx, 2).Value = $chg[1]
    $fieldsSheet.Cells.Item($rowIdx, 4).Value = $chg[2]
}

Write-Host "Added classNumberOfLines and methodNumberOfLines sheets; reordered classFields rows."
